# "update scripts wuth new tpm" -- re-run of the NATMI TPM pipeline produced new
# ligand/receptor statistics for Cd274-Pdcd1, and also expanded the result table
# from 5 rows (one per Sending cluster, Target cluster fixed to Inflammatory-Mac)
# to 10 rows (every Sending cluster x {Inflammatory-Mac, Resolving-Mac} Target cluster).
# Columns (A:T) are unchanged:
#   A: Sending cluster
#   B: Ligand symbol
#   C: Receptor symbol
#   D: Target cluster
#   E: Ligand-expressing cells
#   F: Ligand detection rate
#   G: Ligand average expression value
#   H: Ligand total expression value
#   I: Ligand derived specificity of average expression value
#   J: Ligand derived specificity of total expression value
#   K: Receptor-expressing cells
#   L: Receptor detection rate
#   M: Receptor average expression value
#   N: Receptor total expression value
#   O: Receptor derived specificity of average expression value
#   P: Receptor derived specificity of total expression value
#   Q: Edge average expression weight
#   R: Edge total expression weight
#   S: Edge average expression derived specificity
#   T: Edge total expression derived specificity

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Inflammatory-Mac (Cd274-Pdcd1)
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Cd274"
$ws.Cells.Item(2, 3).Value = "Pdcd1"
$ws.Cells.Item(2, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 12.780993
$ws.Cells.Item(2, 8).Value = 38.342979
$ws.Cells.Item(2, 9).Value = 0.1912157377894449
$ws.Cells.Item(2, 10).Value = 0.1923437725816443
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.6671056666666667
$ws.Cells.Item(2, 14).Value = 2.001317
$ws.Cells.Item(2, 15).Value = 0.91238273339576
$ws.Cells.Item(2, 16).Value = 0.91238273339576
$ws.Cells.Item(2, 17).Value = 8.526272855927001
$ws.Cells.Item(2, 18).Value = 76.73645570334301
$ws.Cells.Item(2, 19).Value = 0.1744619375126207
$ws.Cells.Item(2, 20).Value = 0.1754911369796931

# Row 3: ECs -> Resolving-Mac (Cd274-Pdcd1)
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Cd274"
$ws.Cells.Item(3, 3).Value = "Pdcd1"
$ws.Cells.Item(3, 4).Value = "Resolving-Mac"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 12.780993
$ws.Cells.Item(3, 8).Value = 38.342979
$ws.Cells.Item(3, 9).Value = 0.1912157377894449
$ws.Cells.Item(3, 10).Value = 0.1923437725816443
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.064063
$ws.Cells.Item(3, 14).Value = 0.192189
$ws.Cells.Item(3, 15).Value = 0.08761726660423996
$ws.Cells.Item(3, 16).Value = 0.08761726660423996
$ws.Cells.Item(3, 17).Value = 0.818788754559
$ws.Cells.Item(3, 18).Value = 7.369098791031
$ws.Cells.Item(3, 19).Value = 0.01675380027682424
$ws.Cells.Item(3, 20).Value = 0.01685263560195123

# Row 4: FAPs -> Inflammatory-Mac (Cd274-Pdcd1)
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Cd274"
$ws.Cells.Item(4, 3).Value = "Pdcd1"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.763789333333333
$ws.Cells.Item(4, 8).Value = 8.291368
$ws.Cells.Item(4, 9).Value = 0.04134890117441825
$ws.Cells.Item(4, 10).Value = 0.04159282983679289
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.6671056666666667
$ws.Cells.Item(4, 14).Value = 2.001317
$ws.Cells.Item(4, 15).Value = 0.91238273339576
$ws.Cells.Item(4, 16).Value = 0.91238273339576
$ws.Cells.Item(4, 17).Value = 1.843739525739556
$ws.Cells.Item(4, 18).Value = 16.593655731656
$ws.Cells.Item(4, 19).Value = 0.03772602347642687
$ws.Cells.Item(4, 20).Value = 0.03794857977615782

# Row 5: FAPs -> Resolving-Mac (Cd274-Pdcd1)
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Cd274"
$ws.Cells.Item(5, 3).Value = "Pdcd1"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.763789333333333
$ws.Cells.Item(5, 8).Value = 8.291368
$ws.Cells.Item(5, 9).Value = 0.04134890117441825
$ws.Cells.Item(5, 10).Value = 0.04159282983679289
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.064063
$ws.Cells.Item(5, 14).Value = 0.192189
$ws.Cells.Item(5, 15).Value = 0.08761726660423996
$ws.Cells.Item(5, 16).Value = 0.08761726660423996
$ws.Cells.Item(5, 17).Value = 0.1770566360613333
$ws.Cells.Item(5, 18).Value = 1.593509724552
$ws.Cells.Item(5, 19).Value = 0.003622877697991375
$ws.Cells.Item(5, 20).Value = 0.003644250060635069

# Row 6: Inflammatory-Mac -> Inflammatory-Mac (Cd274-Pdcd1)
$ws.Cells.Item(6, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(6, 2).Value = "Cd274"
$ws.Cells.Item(6, 3).Value = "Pdcd1"
$ws.Cells.Item(6, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 29.29107666666667
$ws.Cells.Item(6, 8).Value = 87.87323
$ws.Cells.Item(6, 9).Value = 0.4382221972474175
$ws.Cells.Item(6, 10).Value = 0.4408073918078855
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.6671056666666667
$ws.Cells.Item(6, 14).Value = 2.001317
$ws.Cells.Item(6, 15).Value = 0.91238273339576
$ws.Cells.Item(6, 16).Value = 0.91238273339576
$ws.Cells.Item(6, 17).Value = 19.54024322710112
$ws.Cells.Item(6, 18).Value = 175.86218904391
$ws.Cells.Item(6, 19).Value = 0.3998263661592946
$ws.Cells.Item(6, 20).Value = 0.4021850530387343

# Row 7: Inflammatory-Mac -> Resolving-Mac (Cd274-Pdcd1)
$ws.Cells.Item(7, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(7, 2).Value = "Cd274"
$ws.Cells.Item(7, 3).Value = "Pdcd1"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 29.29107666666667
$ws.Cells.Item(7, 8).Value = 87.87323
$ws.Cells.Item(7, 9).Value = 0.4382221972474175
$ws.Cells.Item(7, 10).Value = 0.4408073918078855
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.064063
$ws.Cells.Item(7, 14).Value = 0.192189
$ws.Cells.Item(7, 15).Value = 0.08761726660423996
$ws.Cells.Item(7, 16).Value = 0.08761726660423996
$ws.Cells.Item(7, 17).Value = 1.876474244496667
$ws.Cells.Item(7, 18).Value = 16.88826820047
$ws.Cells.Item(7, 19).Value = 0.03839583108812281
$ws.Cells.Item(7, 20).Value = 0.03862233876915117

# Row 8: MuSCs -> Inflammatory-Mac (Cd274-Pdcd1)
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Cd274"
$ws.Cells.Item(8, 3).Value = "Pdcd1"
$ws.Cells.Item(8, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.5
$ws.Cells.Item(8, 7).Value = 1.175998
$ws.Cells.Item(8, 8).Value = 2.351996
$ws.Cells.Item(8, 9).Value = 0.01759404180949881
$ws.Cells.Item(8, 10).Value = 0.01179855596866736
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.6671056666666667
$ws.Cells.Item(8, 14).Value = 2.001317
$ws.Cells.Item(8, 15).Value = 0.91238273339576
$ws.Cells.Item(8, 16).Value = 0.91238273339576
$ws.Cells.Item(8, 17).Value = 0.7845149297886668
$ws.Cells.Item(8, 18).Value = 4.707089578732001
$ws.Cells.Item(8, 19).Value = 0.01605249995762981
$ws.Cells.Item(8, 20).Value = 0.01076479874481558

# Row 9: MuSCs -> Resolving-Mac (Cd274-Pdcd1)
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Cd274"
$ws.Cells.Item(9, 3).Value = "Pdcd1"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.5
$ws.Cells.Item(9, 7).Value = 1.175998
$ws.Cells.Item(9, 8).Value = 2.351996
$ws.Cells.Item(9, 9).Value = 0.01759404180949881
$ws.Cells.Item(9, 10).Value = 0.01179855596866736
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.064063
$ws.Cells.Item(9, 14).Value = 0.192189
$ws.Cells.Item(9, 15).Value = 0.08761726660423996
$ws.Cells.Item(9, 16).Value = 0.08761726660423996
$ws.Cells.Item(9, 17).Value = 0.07533795987400001
$ws.Cells.Item(9, 18).Value = 0.452027759244
$ws.Cells.Item(9, 19).Value = 0.001541541851869002
$ws.Cells.Item(9, 20).Value = 0.001033757223851775

# Row 10: Resolving-Mac -> Inflammatory-Mac (Cd274-Pdcd1)
$ws.Cells.Item(10, 1).Value = "Resolving-Mac"
$ws.Cells.Item(10, 2).Value = "Cd274"
$ws.Cells.Item(10, 3).Value = "Pdcd1"
$ws.Cells.Item(10, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 20.828839
$ws.Cells.Item(10, 8).Value = 62.48651700000001
$ws.Cells.Item(10, 9).Value = 0.3116191219792205
$ws.Cells.Item(10, 10).Value = 0.31345744980501
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.6671056666666667
$ws.Cells.Item(10, 14).Value = 2.001317
$ws.Cells.Item(10, 15).Value = 0.91238273339576
$ws.Cells.Item(10, 16).Value = 0.91238273339576
$ws.Cells.Item(10, 17).Value = 13.89503652698767
$ws.Cells.Item(10, 18).Value = 125.055328742889
$ws.Cells.Item(10, 19).Value = 0.284315906289788
$ws.Cells.Item(10, 20).Value = 0.2859931648563593

# Row 11: Resolving-Mac -> Resolving-Mac (Cd274-Pdcd1)
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Cd274"
$ws.Cells.Item(11, 3).Value = "Pdcd1"
$ws.Cells.Item(11, 4).Value = "Resolving-Mac"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 20.828839
$ws.Cells.Item(11, 8).Value = 62.48651700000001
$ws.Cells.Item(11, 9).Value = 0.3116191219792205
$ws.Cells.Item(11, 10).Value = 0.31345744980501
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.064063
$ws.Cells.Item(11, 14).Value = 0.192189
$ws.Cells.Item(11, 15).Value = 0.08761726660423996
$ws.Cells.Item(11, 16).Value = 0.08761726660423996
$ws.Cells.Item(11, 17).Value = 1.334357912857
$ws.Cells.Item(11, 18).Value = 12.009221215713
$ws.Cells.Item(11, 19).Value = 0.02730321568943254
$ws.Cells.Item(11, 20).Value = 0.02746428494865073
